# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Re-order two country-name pairs in the country list (values, not the
#    shared-string table, since the engine rebuilds/compacts sst on save)
#  - Bump the "Datos actualizados" timestamp string
#  - Refresh a handful of per-country case counters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 08:22"

# --- Alemania (row 8): recuperados/activos swap ------------------------
$ws.Range("D8").Value = 64300
$ws.Range("E8").Value = 60532

# --- Bielorrusia / Tailandia swap positions (rows 51-52) ---------------
$ws.Range("A51").Value = "Tailandia"
$ws.Range("B51").Value = 2579
$ws.Range("C51").Value = 28
$ws.Range("D51").Value = 1288
$ws.Range("E51").Value = 1251
$ws.Range("F51").Value = 61
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 40

$ws.Range("A52").Value = "Bielorrusia"
$ws.Range("B52").Value = 2578
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 203
$ws.Range("E52").Value = 2349
$ws.Range("F52").Value = 50
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 26

# --- Moldavia (row 60) ---------------------------------------------------
$ws.Range("E60").Value = 1536
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 32

# --- Hungria (row 63) -----------------------------------------------------
$ws.Range("B63").Value = 1458
$ws.Range("C63").Value = 48
$ws.Range("D63").Value = 120
$ws.Range("E63").Value = 1229
$ws.Range("G63").Value = 10
$ws.Range("H63").Value = 109

# --- Kazajistan (row 75) ---------------------------------------------------
$ws.Range("B75").Value = 973
$ws.Range("C75").Value = 22
$ws.Range("E75").Value = 863
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 11

# --- Taiwan / Reunion / Jordania reorder (rows 100-102) ---------------------
$ws.Range("A100").Value = "Taiwan"
$ws.Range("B100").Value = 393
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 114
$ws.Range("E100").Value = 273
$ws.Range("F100").Value = 0
$ws.Range("H100").Value = 6

$ws.Range("A101").Value = "Reunion"
$ws.Range("D101").Value = 40
$ws.Range("E101").Value = 349
$ws.Range("F101").Value = 3
$ws.Range("H101").Value = 0

$ws.Range("A102").Value = "Jordania"
$ws.Range("B102").Value = 389
$ws.Range("D102").Value = 201
$ws.Range("E102").Value = 181
$ws.Range("F102").Value = 5
$ws.Range("H102").Value = 7

# --- El Salvador (row 126) ---------------------------------------------------
$ws.Range("F126").Value = 2
